# V 0.51-B48 Doc update
# Insert a new "DEPARR" flag column into Tabelle2, just before the old
# END_OF_COL marker column (i.e. before the previous "EC" column), shifting
# the trailing END_OF_COL / Title columns one position to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle2")

# Insert a new blank column at EC (pushes old EC->ED, old ED->EE, etc.)
$ws.Columns("EC:EC").Insert()

# Header row: new column name
$ws.Range("EC1").Value = "DEPARR"

# Data rows: same filler value ("|") as the other flag columns in this block
$ws.Range("EC2:EC40").Value = "|"

# Leave the active selection where the author ended up after the edit
$ws.Range("CZ14").Select()
